# Weekly data refresh: insert a new price record as row 416, pushing the
# existing rows 416:443 down to 417:444 (dimension grows from A1:T443 to
# A1:T444).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 416; Excel shifts 416:443 down to 417:444.
$ws.Rows.Item(416).Insert()

# Populate the newly inserted row 416 with the new weekly record.
$ws.Range("A416").Value2 = 9
$ws.Range("B416").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C416").Value2 = "Metropolitana"
$ws.Range("D416").Value2 = 44746
$ws.Range("D416").NumberFormat = $ws.Range("D417").NumberFormat
$ws.Range("E416").Value2 = 13
$ws.Range("F416").Value2 = "Fruta"
$ws.Range("G416").Value2 = 100108
$ws.Range("H416").Value2 = "Tropicales y subtropicales"
$ws.Range("I416").Value2 = 100108002
$ws.Range("J416").Value2 = "Mango"
$ws.Range("K416").Value2 = "Sin especificar"
$ws.Range("L416").Value2 = "Primera"
$ws.Range("M416").Value2 = 420
$ws.Range("N416").Value2 = 7500
$ws.Range("O416").Value2 = 8000
$ws.Range("P416").Value2 = 7738
$ws.Range("Q416").Value2 = "$/bandeja 4 kilos"
$ws.Range("R416").Value2 = "Brasil"
$ws.Range("S416").Value2 = 1934
$ws.Range("T416").Value2 = 4
